$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.498.28'
$ws.Range('E2').Value = '  +1.32%  '

$ws.Range('D3').Value = '3.025.63'
$ws.Range('E3').Value = '  -0.11%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').Value = '547.12'
$ws.Range('E5').Value = '  +1.79%  '

$ws.Range('D6').Value = '136.18'
$ws.Range('E6').Value = '  +0.98%  '

$ws.Range('E7').Value = '  -0.15%  '

$ws.Range('D8').Value = '3.022.60'
$ws.Range('E8').Value = '  +0.11%  '

$ws.Range('D9').Value = '0.495'
$ws.Range('E9').Value = '  -0.82%  '

$ws.Range('D10').Value = '6.22'
$ws.Range('E10').Value = '  +1.09%  '

$ws.Range('D11').Value = '0.147'
$ws.Range('E11').Value = '  -1.90%  '

$ws.Range('D12').Value = '0.445'
$ws.Range('E12').Value = '  -1.47%  '

$ws.Range('E13').Value = '  +0.30%  '

$ws.Range('D14').Value = '34.20'
$ws.Range('E14').Value = '  -1.74%  '

$ws.Range('D15').Value = '3.531.04'
$ws.Range('E15').Value = '  +0.46%  '

$ws.Range('D16').Value = '62.586.50'
$ws.Range('E16').Value = '  +1.43%  '

$ws.Range('D17').Value = '3.032.37'
$ws.Range('E17').Value = '  +0.21%  '

$ws.Range('D18').Value = '0.108'
$ws.Range('E18').Value = '  -2.42%  '

$ws.Range('D19').Value = '6.63'
$ws.Range('E19').Value = '  -0.74%  '

$ws.Range('D20').Value = '475.67'
$ws.Range('E20').Value = '  +1.29%  '

$ws.Range('D21').Value = '13.44'
$ws.Range('E21').Value = '  +0.35%  '

$ws.Range('D22').Value = '0.663'
$ws.Range('E22').Value = '  -3.20%  '

$ws.Range('D23').Value = '7.08'
$ws.Range('E23').Value = '  +0.50%  '

$ws.Range('D24').Value = '79.95'
$ws.Range('E24').Value = '  +0.24%  '

$ws.Range('D25').Value = '12.23'
$ws.Range('E25').Value = '  -0.18%  '

$ws.Range('E26').Value = '  -0.29%  '

$ws.Range('D27').Value = '2.72'
$ws.Range('E27').Value = '  +0.36%  '

$ws.Range('D28').Value = '7.74'
$ws.Range('E28').Value = '  -2.93%  '

$ws.Range('D30').Value = '1.95'
$ws.Range('E30').Value = '  +1.82%  '

$ws.Range('D31').Value = '25.64'
$ws.Range('E31').Value = '  -0.57%  '

$ws.Range('E32').Value = '  -1.20%  '

$ws.Range('D33').Value = '2.37'
$ws.Range('E33').Value = '  +2.24%  '

$ws.Range('D34').Value = '5.59'
$ws.Range('E34').Value = '  +0.21%  '

$ws.Range('D35').Value = '54.63'
$ws.Range('E35').Value = '  -1.72%  '

$ws.Range('D36').Value = '5.87'
$ws.Range('E36').Value = '  -1.79%  '

$ws.Range('D37').Value = '456.21'
$ws.Range('E37').Value = '  -3.37%  '

$ws.Range('D38').Value = '0.0803'
$ws.Range('E38').Value = '  +0.52%  '

$ws.Range('D39').Value = '3.062.30'
$ws.Range('E39').Value = '  -4.66%  '

$ws.Range('D40').Value = '0.0388'
$ws.Range('E40').Value = '  -1.21%  '

$ws.Range('E41').Value = '  -2.18%  '

$ws.Range('D42').Value = '8.14'
$ws.Range('E42').Value = '  -0.94%  '

$ws.Range('D43').Value = '2.48'
$ws.Range('E43').Value = '  -1.53%  '

$ws.Range('D44').Value = '27.27'
$ws.Range('E44').Value = '  +0.57%  '

$ws.Range('E45').Value = '  -0.13%  '

$ws.Range('D46').Value = '0.247'
$ws.Range('E46').Value = '  -2.09%  '

$ws.Range('D47').Value = '2.00'
$ws.Range('E47').Value = '  -1.78%  '

$ws.Range('E48').Value = '  -1.14%  '

$ws.Range('D49').Value = '116.06'
$ws.Range('E49').Value = '  -3.21%  '

$ws.Range('D50').Value = '0.0₃0496'
$ws.Range('E50').Value = '  -1.56%  '

$ws.Range('B51').Value = 'BitgetToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
$ws.Range('D51').Value = '1.26'
$ws.Range('E51').Value = '  -0.32%  '
